# Daily attendance processing - 2025-12-21 07:28:34
# Normalises the "Recorded By" (column G) cell values: the literal
# "System" token that was previously listed first is reordered so the
# human/automation identifier is listed first instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Exact before -> after replacements for the "Recorded By" column (G).
$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
